# Auto-generated edit script: update quotation style (" -> ') in specific
# English (en_US, column C) story lines per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-LeadingApostropheText($row, $col, $formulaText) {
    # Excel's normal Value-setter strips a literal leading apostrophe
    # (it's interpreted as the text-prefix marker), so build the text via a
    # formula, compute it, then paste back the computed value only.
    $c = $ws.Cells.Item($row, $col)
    $c.Formula = $formulaText
    $c.Copy() | Out-Null
    $c.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

$ws.Cells.Item(15, 3).Value = "It needs an 'and then' in there.`n"
$f20 = "=CHAR(39)&""District 14 reporting, all clear.""&CHAR(39)&CHAR(10)"
Set-LeadingApostropheText 20 3 $f20
$f21 = "=CHAR(39)&""District 5 reporting, situation normal.""&CHAR(39)&CHAR(10)"
Set-LeadingApostropheText 21 3 $f21
$f22 = "=CHAR(39)&""District 1 reporting, nothing out of the ordinary.""&CHAR(39)&CHAR(10)"
Set-LeadingApostropheText 22 3 $f22
$f23 = "=CHAR(39)&""Green across the board. Good work, everyone. Stay on your toes.""&CHAR(39)&CHAR(10)"
Set-LeadingApostropheText 23 3 $f23
$f24 = "=CHAR(39)&""Roger.""&CHAR(39)&CHAR(10)"
Set-LeadingApostropheText 24 3 $f24
$f25 = "=CHAR(39)&""Understood.""&CHAR(39)&CHAR(10)"
Set-LeadingApostropheText 25 3 $f25
$ws.Cells.Item(30, 3).Value = "[name=""Serious L.G.D. Officer""]   ...'Suspected Originium explosive device activation'? What is it? Somebody shooting off fireworks in the park?`n"
$ws.Cells.Item(37, 3).Value = "[name=""Frivolous L.G.D. Officer""]   I've never seen the 'Nian' before. What is it?`n"
$ws.Cells.Item(77, 3).Value = "[name=""Snowsant""]   My grandma told me stories about the 'Nian' when I was little... But they were so scary I never thought they could be true...`n"
$ws.Cells.Item(156, 3).Value = "[name=""Ch'en""]   Obstruction of justice, destruction of public property, trespassing, disturbing the peace, assault. You made a lot of noise last year, 'Madame Oni.'`n"
$ws.Cells.Item(175, 3).Value = "[name=""Madame Oni""]   'This year's special circumstances...' You know what I mean, Ch'en Sir?`n"
$ws.Cells.Item(240, 3).Value = "[Decision(options=""To 'ring in the New Year,' right?;What's this 'Nian' anyway?;......"",values=""1;2;3"")]`n"
$ws.Cells.Item(241, 3).Value = "[name=""Amiya""]   The 'Nian' is a legend... or so they say.`n"
$ws.Cells.Item(247, 3).Value = "[name=""Amiya""]   There's even a conspiracy theory floating around that some organized crime syndicate or other shadowy group fabricated the whole 'Nian' thing to cover up their activities...`n"
$ws.Cells.Item(260, 3).Value = "[name=""Madame Oni""]   What ngaang gaau 'get-togethers' are you talking about... bullshit... ugh...`n"
$ws.Cells.Item(274, 3).Value = "[name=""Ch'en""]   When she wakes up, you tell her, the L.G.D. got a lead on the 'Nian.'`n"
$ws.Cells.Item(306, 3).Value = "[name=""Serious L.G.D. Officer""]   ...'Tales of the Nian: A Compilation of Catastrophes from Several Countries'?`n"
